# Change the table style ("Table_0" -> the built-in PowerPoint table style
# with GUID {2C9CB8B2-B285-462C-ACB7-E517F99726C5}) on the table that lives
# on slide 6 of the deck.
#
# PowerPoint's Table object doesn't allow StyleId to be assigned directly
# (it throws "Table styles cannot be assigned through a property - call
# Table.ApplyStyle(...) instead"), so we use Table.ApplyStyle().

$p = $ppt.ActivePresentation

$targetStyleId = "{2C9CB8B2-B285-462C-ACB7-E517F99726C5}"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($targetStyleId)
        }
    }
}
